$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (support Zone) with new ticker symbols for rows 2-8
$ws.Range("C2").Value = "NSE:AARTIPHARM"
$ws.Range("C3").Value = "NSE:CARERATING"
$ws.Range("C4").Value = "NSE:GOCOLORS"
$ws.Range("C5").Value = "NSE:HUBTOWN"
$ws.Range("C6").Value = "NSE:JUBLPHARMA"
$ws.Range("C7").Value = "NSE:PGIL"
$ws.Range("C8").Value = "NSE:PPAP"

# Clear column B (Buying Opportunity) for rows 2-5 (B6-B8 already empty)
$ws.Range("B2:B5").ClearContents()

# Clear column E (Short buildup) for rows 2-8
$ws.Range("E2:E8").ClearContents()

# Delete rows 9-16 entirely (shrinking the used range to A1:F8)
$ws.Range("A9:A16").EntireRow.Delete()
